$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the points list for Cowlevel001 (row 19, column C)
$ws.Cells.Item(19, 3).Value = "506276.25|42784.48|46083.01||495329.25|43252.48|46083.01||497870.25|53525.48|46083.01||497870.25|53525.48|47007||503903.25|53525.48|47007||504647.25|47997.48|47007||499929.25|47997.48|47007"

# Row 19 height changes from 49.5 to 33 (wrap text shrinks because the C value is shorter now)
$ws.Rows.Item(19).RowHeight = 33

# Fill in the points column (C) for the existing Cowlevel002..Cowlevel006 rows (20..24)
$ws.Cells.Item(20, 3).Value = "598563.75|594154.88|7143.97||602921.75|595009.88|7143.97||602921.75|595009.88|7851.97||602921.75|601241.88|7851.97||598241.75|602673.88|7851.97||598241.75|594819.88|7851.97||602921.75|600363.88|7295.97"
$ws.Cells.Item(21, 3).Value = "498910.44|46868.74|45986||496482.44|48743.74|45986||498294.44|55570.74|45986||499992.44|56135.74|45986||499992.44|56135.74|46753.00||499228.44|47348.74|46753||496626.44|48303.74|46753||498189.44|55497.74|46424"
$ws.Cells.Item(22, 3).Value = "97705.11|1001785.63|6204.8||97705.11|995826.63|6204.8||101330.11|995826.63|6204.8||101330.11|1001359.63|6204.80||101330.11|1001359.63|6911.80||97885.11|1001359.63|6911.80||97885.11|996353.63|6911.80||101778.11|996223.63|6911.80"
$ws.Cells.Item(23, 3).Value = "98879.05|1000148.06|10082||98879.05|997982.06|10082||100784.05|997982.06|10082||100186.05|999942.06|10082||99185.45|999750.75|10589||99239.05|998549.06|10590||99851.05|998549.06|10590||99851.05|999709.06|10590"
$ws.Cells.Item(24, 3).Value = "1010497.13|18739.56|4593.87||1006911.13|16723.56|4593.87||1004790.13|19763.56|4593.87||1006629.13|19763.56|5473.87||1009839.13|18905.56|5473.87||1007139.13|17599.56|5473.87"

# Add the new Cowlevel007 row (25)
$ws.Cells.Item(25, 1).Value = 21
$ws.Cells.Item(25, 2).Value = "Cowlevel007"
$ws.Cells.Item(25, 3).Value = "8978.4|229.39|296699.97||5201.4|229.39|296699.97||5201.4|5546.39|296699.97||6494.4|5546.39|296699.97||6494.4|5037.39|297305.97||6494.4|1696.39|297305.97||5428.4|3174.39|297305.97"

# Update sheet view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C21").Select()
